$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 674, shifting existing rows
# 674-681 down to 677-684.
$null = $ws.Rows("674:676").Insert()

# Populate the 3 new rows with the new order data (remessa 80267652).
# Column A holds a purely numeric-looking value ("80267652") that must be
# stored as text (matching the rest of the "Remessa" column). Writing it
# through a text formula and then pasting the calculated value back avoids
# Excel auto-marking the cell with a quote-prefix / text number format
# (which would otherwise register an extra, unused cell style).
$ws.Range("A674").Formula = '="80267652"'
$ws.Range("A675").Formula = '="80267652"'
$ws.Range("A676").Formula = '="80267652"'
$ws.Range("A674:A676").Copy()
$null = $ws.Range("A674:A676").PasteSpecial(-4163)

$ws.Range("B674").Value = "00041-DIG-I"
$ws.Range("C674").Value = 500

$ws.Range("B675").Value = "13116-VIS-I"
$ws.Range("C675").Value = 50

$ws.Range("B676").Value = "15663-TDK-N"
$ws.Range("C676").Value = 420

# Re-apply the formatting used by the surrounding data rows (columns A/B
# use style index 3) so the inserted rows match the rest of the sheet.
$ws.Range("A673").Copy()
$null = $ws.Range("A674:A676").PasteSpecial(-4122)
$ws.Range("B673").Copy()
$null = $ws.Range("B674:B676").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to reflect where the author last clicked.
$null = $ws.Range("H12").Select()
